# Portal Check Added and Asana Updates
#
# The SKU/pricing detail rows for the "QVR" (Quantity Variance) transaction
# got re-matched: the line that used to be reported on row 2 (SKU 100284)
# actually belongs on row 4, and vice-versa for the SKU 1572435 line. Row 3
# (SKU 24531799) is untouched. This swaps the SKU / Vendor Part # / Units
# Received / Units Invoiced / PO Unit Cost / Invoice Unit Cost / Extended
# Cost Variance values between row 2 and row 4.
#
# Values in this sheet are stored as text (shared strings) even though many
# look numeric (e.g. "124", "15.9", "-63.6"). Plain `Range.Value = "124"`
# would be auto-coerced to a number by the host, so each value is written
# via a temporary text formula and then frozen into a literal with
# Copy + PasteSpecial(values) - this keeps the cell's stored type as text
# and avoids introducing any new cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param(
        [string]$Address,
        [string]$Text
    )

    $escaped = $Text.Replace('"', '""')
    $ws.Range($Address).Formula = '="' + $escaped + '"'
    $ws.Range($Address).Copy()
    $ws.Range($Address).PasteSpecial(-4163)
}

# --- Row 2 gets what used to be row 4's data (SKU 1572435 line) ---
Set-TextValue "I2" "1572435"
Set-TextValue "J2" "920-006481"
Set-TextValue "L2" "60"
Set-TextValue "M2" "64"
Set-TextValue "N2" "31.43"
Set-TextValue "O2" "31.43"
Set-TextValue "P2" "-125.72"

# --- Row 4 gets what used to be row 2's data (SKU 100284 line) ---
Set-TextValue "I4" "100284"
Set-TextValue "J4" "981-000507"
Set-TextValue "L4" "124"
Set-TextValue "M4" "128"
Set-TextValue "N4" "15.9"
Set-TextValue "O4" "15.9"
Set-TextValue "P4" "-63.6"
